$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.92"
$ws.Range("G2").Value = "'10"
$ws.Range("D3").Value = "'21.95"
$ws.Range("G3").Value = "'10"
$ws.Range("D4").Value = "'5.395"
$ws.Range("G4").Value = "'10"
$ws.Range("D5").Value = "'0.05990"
$ws.Range("G5").Value = "'10"
$ws.Range("G6").Value = "'10"
$ws.Range("D7").Value = "'6.401"
$ws.Range("G7").Value = "'10"
$ws.Range("D8").Value = "'0.8104"
$ws.Range("G8").Value = "'10"
$ws.Range("D9").Value = "'0.9536"
$ws.Range("G9").Value = "'10"
$ws.Range("D10").Value = "'0.1431"
$ws.Range("G10").Value = "'10"
$ws.Range("D11").Value = "'0.07386"
$ws.Range("G11").Value = "'10"
$ws.Range("G12").Value = "'10"
$ws.Range("D13").Value = "'0.03054"
$ws.Range("G13").Value = "'10"
$ws.Range("D14").Value = "'0.09420"
$ws.Range("G14").Value = "'10"
$ws.Range("D15").Value = "'3.999"
$ws.Range("G15").Value = "'10"
$ws.Range("D16").Value = "'0.001589"
$ws.Range("G16").Value = "'10"
$ws.Range("D17").Value = "'0.04800"
$ws.Range("G17").Value = "'10"
$ws.Range("B18").Value = "'One"
$ws.Range("C18").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005871"
$ws.Range("E18").Value = "'17OneONEWorstin24h"
$ws.Range("G18").Value = "'10"
$ws.Range("B19").Value = "'TigerCash"
$ws.Range("C19").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006132"
$ws.Range("E19").Value = "'18TigerCashTCH"
$ws.Range("G19").Value = "'10"
$ws.Range("B20").Value = "'HotbitToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.005069"
$ws.Range("E20").Value = "'19HotbitTokenHTB"
$ws.Range("G20").Value = "'10"
$ws.Range("B21").Value = "'BitKan"
$ws.Range("C21").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009868"
$ws.Range("E21").Value = "'20BitKanKAN"
$ws.Range("G21").Value = "'10"
$ws.Range("B22").Value = "'NitroEx"
$ws.Range("C22").Value = "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001000"
$ws.Range("E22").Value = "'21NitroExNTX"
$ws.Range("G22").Value = "'10"
$ws.Range("B23").Value = "'LEO"
$ws.Range("C23").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.704"
$ws.Range("E23").Value = "'22LEOLEO"
$ws.Range("G23").Value = "'10"
$ws.Range("B24").Value = "'BTSEToken"
$ws.Range("C24").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.185"
$ws.Range("E24").Value = "'23BTSETokenBTSE"
$ws.Range("G24").Value = "'10"
$ws.Range("G25").Value = "'10"
$ws.Range("D26").Value = "'0.1284"
$ws.Range("G26").Value = "'10"
$ws.Range("G27").Value = "'10"
$ws.Range("G28").Value = "'10"
$ws.Range("G29").Value = "'10"
$ws.Range("G30").Value = "'10"
$ws.Range("G31").Value = "'10"
$ws.Range("G32").Value = "'10"
$ws.Range("G33").Value = "'10"
$ws.Range("G34").Value = "'10"
$ws.Range("G35").Value = "'10"
$ws.Range("G36").Value = "'10"
$ws.Range("G37").Value = "'10"
$ws.Range("G38").Value = "'10"
$ws.Range("G39").Value = "'10"
$ws.Range("D40").Value = "'0.04004"
$ws.Range("G40").Value = "'10"
$ws.Range("D41").Value = "'0.006523"
$ws.Range("G41").Value = "'10"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("G42").Value = "'10"
$ws.Range("G43").Value = "'10"
$ws.Range("D44").Value = "'0.005311"
$ws.Range("G44").Value = "'10"
$ws.Range("D45").Value = "'0.00005253"
$ws.Range("G45").Value = "'10"
$ws.Range("G46").Value = "'10"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOINBestin24h"
$ws.Range("G47").Value = "'10"
$ws.Range("D48").Value = "'0.02412"
$ws.Range("E48").Value = "'47BOLOBOLO"
$ws.Range("G48").Value = "'10"
$ws.Range("G49").Value = "'10"
$ws.Range("G50").Value = "'10"
$ws.Range("G51").Value = "'10"
